$d = $word.ActiveDocument

$pairs = @(
    @("614×6=", "259×5="),
    @("713×3=", "564×2="),
    @("133×8=", "516×3="),
    @("342×2=", "862×3="),
    @("781×9=", "426×7="),
    @("124×7=", "646×6="),
    @("723×3=", "449×3="),
    @("309×6=", "610×8="),
    @("249×3=", "780×2="),
    @("898×3=", "722×9="),
    @("749×5=", "388×5="),
    @("271×8=", "526×7="),
    @("664×6=", "981×8="),
    @("884×3=", "693×9="),
    @("404×2=", "936×5="),
    @("540×7=", "701×6="),
    @("801×3=", "388×3="),
    @("796×6=", "122×9="),
    @("935×2=", "780×6="),
    @("498×2=", "804×7="),
    @("301×4=", "917×7="),
    @("164×9=", "419×2="),
    @("838×8=", "141×5="),
    @("671×9=", "963×3="),
    @("350×7=", "843×5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
